$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText {
    param($addr, $text)
    # Force the literal text (e.g. "01-08-2022") into the cell without
    # Excel's automatic date-recognition turning it into a date serial
    # number. Pre-marking the cell as Text keeps the typed value literal;
    # ClearFormats afterwards drops the temporary Text format so the cell
    # is left with its original (default/general) style.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Update date strings (slash -> dash format) for rows 3-21
Set-DateText "A3"  "28-07-2022"
Set-DateText "A4"  "01-08-2022"
Set-DateText "A5"  "04-08-2022"
Set-DateText "A6"  "08-08-2022"
Set-DateText "A7"  "11-08-2022"
Set-DateText "A8"  "15-08-2022"
Set-DateText "A9"  "18-08-2022"
Set-DateText "A10" "22-08-2022"
Set-DateText "A11" "25-08-2022"
Set-DateText "A12" "29-08-2022"
Set-DateText "A13" "01-09-2022"
Set-DateText "A14" "05-09-2022"
Set-DateText "A15" "08-09-2022"
Set-DateText "A16" "12-09-2022"
Set-DateText "A17" "15-09-2022"
Set-DateText "A18" "19-09-2022"
Set-DateText "A19" "22-09-2022"
Set-DateText "A20" "26-09-2022"
Set-DateText "A21" "29-09-2022"

# Update attendance counts for affected rows
# Row 3: D 0->1, G 0->1
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4: D 0->1, E 0->1, H 1->0
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5: D 0->1, E 0->1, H 1->0
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 10: D 0->1, E 0->1, H 1->0
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("H10").Value = 0
